# Insert a new data row at row 84 (pushing existing rows 84-132 down to 85-133)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84 and below down by one row, then fill in the new row 84 data.
$ws.Rows.Item(84).Insert()

$ws.Cells.Item(84, 1).Value = 11
$ws.Cells.Item(84, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(84, 3).Value = "Bíobío"
$ws.Cells.Item(84, 4).Value = 44469
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(84, 6).Value = 100112021
$ws.Cells.Item(84, 7).Value = "Ají"
$ws.Cells.Item(84, 8).Value = "Americana (o)"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 40
$ws.Cells.Item(84, 11).Value = 96000
$ws.Cells.Item(84, 12).Value = 98000
$ws.Cells.Item(84, 13).Value = 97000
$ws.Cells.Item(84, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(84, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(84, 16).Value = 3880
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
